$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" (Changed) date column (C) for rows 2-6 from 2023-11-13 to 2023-11-14
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = (Get-Date -Year 2023 -Month 11 -Day 14 -Hour 0 -Minute 0 -Second 0).Date
}
